$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Means" --- add 5-mile and 10-mile radius columns
$ws1 = $wb.Worksheets.Item("Means")

$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# New data columns mirror the existing C/D/E columns: all #NUM! errors
$ws1.Range("F2:F10").Value = "#NUM!"
$ws1.Range("G2:G10").Value = "#NUM!"

# Updated National Average values (recomputed with the new buffers)
$ws1.Range("B9").Value = 29
$ws1.Range("B10").Value = 0.37

# --- Sheet 2: "Standard Deviations" --- add 5-mile and 10-mile radius columns
$ws2 = $wb.Worksheets.Item("Standard Deviations")

$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# New data columns mirror the existing C/D/E columns: all 0
$ws2.Range("F2:F10").Value = 0
$ws2.Range("G2:G10").Value = 0

# Updated National Average SD values (recomputed with the new buffers)
$ws2.Range("B9").Value = 10
$ws2.Range("B10").Value = 0.14
